$d = $word.ActiveDocument

# Inserts "[" + $fieldName + "]" at the end of $para's range as three
# separate runs (matching the target markup "[", fieldName, "]") by
# temporarily bookmarking the middle token -- this forces a run split at
# the bookmark boundaries -- and then deleting the bookmark again so no
# trace of it remains in the saved document.
function Insert-BracketedField($para, $fieldName) {
    $pStart = $para.Range.Start
    $para.Range.InsertAfter("[" + $fieldName + "]")
    $start = $pStart + 1
    $end = $start + $fieldName.Length
    $inner = $d.Range($start, $end)
    $d.Bookmarks.Add("tmpFieldSplit", $inner) | Out-Null
    $d.Bookmarks.Item("tmpFieldSplit").Delete()
}

# Locate the empty paragraph that immediately follows the "[USER_ADDRESS]"
# paragraph (and immediately precedes the page-break paragraph). This is the
# anchor we build the new KYC_* / COMMITMENT_* paragraphs from.
$anchor = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -eq "[USER_ADDRESS]`r") {
        $anchor = $d.Paragraphs.Item($i + 1)
        break
    }
}

# 1) Fill the existing empty paragraph with [KYC_CITY]
$anchor.Range.InsertAfter("[KYC_CITY]")

# 2) Add a new paragraph with [KYC_STATE]
$anchor.Range.InsertParagraphAfter()
$p = $anchor.Next()
$p.Range.InsertAfter("[KYC_STATE]")

# 3) Add a new paragraph with [KYC_PINCODE]
$p.Range.InsertParagraphAfter()
$p = $p.Next()
$p.Range.InsertAfter("[KYC_PINCODE]")

# 4) Add a new paragraph with [KYC_DOMICILE]
$p.Range.InsertParagraphAfter()
$p = $p.Next()
$p.Range.InsertAfter("[KYC_DOMICILE]")

# 5) Add a new blank paragraph (mirrors the blank paragraph style already
#    used throughout this part of the document)
$p.Range.InsertParagraphAfter()
$p = $p.Next()

# 6) Add a new paragraph containing "[COMMITMENT_TAX_EXEMPT]" split into the
#    3 runs "[", "COMMITMENT_TAX_EXEMPT", "]"
$p.Range.InsertParagraphAfter()
$p = $p.Next()
Insert-BracketedField $p "COMMITMENT_TAX_EXEMPT"

# 7) The paragraph that follows holds the manual page break. Temporarily pull
#    the page-break run out, add the bracketed "[COMMITMENT_COMMITMENT_TYPE]"
#    field (3 runs) into the now-simple paragraph, then restore the page
#    break at the very end -- this keeps the new text ahead of the break,
#    exactly like the target markup.
$p = $p.Next()
$breakStart = $p.Range.End - 2
$breakEnd = $p.Range.End - 1
$breakRange = $d.Range($breakStart, $breakEnd)
$breakRange.Delete()
Insert-BracketedField $p "COMMITMENT_COMMITMENT_TYPE"
$p.Range.InsertAfter([char]12)
